$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new columns ------------------------------------
# Old layout:  A sample_name | B..D means | E..G mean-derived (old) | H..J errs | K..M err-derived (old)
# New layout:  A..D unchanged | E..I mean-derived (new, 5 cols) | J..L errs (unchanged) | M..Q err-derived (new, 5 cols)
# Insert 2 columns before the old "E" block, then 2 more before the (now
# shifted) old "K" block so both derived-quantity blocks grow from 3 to 5
# columns while J..L (err_d15N/17O/18O) land exactly where they should.
$ws.Range("E1:F1").EntireColumn.Insert()
$ws.Range("M1:N1").EntireColumn.Insert()

# --- 2. Header row ----------------------------------------------------
$ws.Range("E1").Value = "mean_d18O/15N"
$ws.Range("F1").Value = "mean_d17O/15N"
$ws.Range("G1").Value = "mean_d15N18O/15N"
$ws.Range("H1").Value = "mean_d17O18O/15N"
$ws.Range("I1").Value = "mean_d18O18O/15N"
$ws.Range("M1").Value = "err_d18O/15N"
$ws.Range("N1").Value = "err_d17O/15N"
$ws.Range("O1").Value = "err_d15N18O/15N"
$ws.Range("P1").Value = "err_d17O18O/15N"
# Q1 already reads "err_d18O18O" after the shift above - no change needed.

# --- 3. Row 2 (USGS34) -------------------------------------------------
$ws.Range("E2").Formula = "=1000*((D2+1000)/(B2+1000)-1)"
$ws.Range("F2").Formula = "=1000*((C2+1000)/(B2+1000)-1)"
# G2/H2/I2 stay blank (no d15N18O/d17O18O/d18O18O data for this standard);
# H2 is additionally flagged with the "needs review" red font + 0.0 format.
$ws.Range("H2").NumberFormat = "0.0"
$ws.Range("H2").Font.Color = 192
$ws.Range("M2").Formula = "=SQRT(J2^2+L2^2)"
$ws.Range("N2").Formula = "=SQRT(K2^2+L2^2)"

# --- 4. Row 3 (USGS35) -------------------------------------------------
$ws.Range("E3").Formula = "=1000*((D3+1000)/(B3+1000)-1)"
$ws.Range("F3").Formula = "=1000*((C3+1000)/(B3+1000)-1)"
$ws.Range("G3").Formula = "=D3"
$ws.Range("H3").Formula = "=((1000+C3)*(1000+E3)/1000)-1000"
$ws.Range("I3").Formula = "=((1000+D3)*(1000+E3)/1000)-1000"
$ws.Range("M3").Formula = "=SQRT(J3^2+L3^2)"
$ws.Range("N3").Formula = "=SQRT(K3^2+L3^2)"
$ws.Range("O3").Formula = "=SQRT(L3^2+(J3)^2)"
$ws.Range("P3").Formula = "=SQRT(J2^2+L2^2)"
$ws.Range("Q3").Formula = "=SQRT(J2^2+L2^2)"

# --- 5. Rows 4-10 (USGS32, NICO1-6) ------------------------------------
for ($r = 4; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=1000*((D$r+1000)/(B$r+1000)-1)"
    $ws.Cells.Item($r, 13).Formula = "=SQRT(J$r^2+L$r^2)"
}

Write-Host "done"
